# Generate Report for Handback
# - Flip status from "Ready for handoff" to "Handed back: in sync with en-US"
#   on the Overview sheet (zh-cn/de-de columns) and on each per-locale sheet.
# - Refresh "Latest Handback DateTime" with the new handback timestamp.
# - Clear "Latest Handback Name" and "Error Detail" now that the handback is
#   in sync (no more stale-version error).
# - Widen/narrow a few columns so the new text fits / the now-empty columns
#   shrink back down.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: status columns for zh-cn (E) and de-de (F), rows 2 and 3
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F3").Value = $statusText
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $statusText
$zhcn.Range("L2").Value = "2017-02-21 04:01:50"
$zhcn.Range("M2").Value = ""
$zhcn.Range("R2").Value = ""

$zhcn.Range("C3").Value = $statusText
$zhcn.Range("L3").Value = "2017-02-21 04:01:50"
$zhcn.Range("M3").Value = ""
$zhcn.Range("R3").Value = ""

$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(13).ColumnWidth = 23.0
$zhcn.Columns.Item(18).ColumnWidth = 12.833333333333332

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $statusText
$dede.Range("L2").Value = "2017-02-21 04:02:14"
$dede.Range("M2").Value = ""
$dede.Range("R2").Value = ""

$dede.Range("C3").Value = $statusText
$dede.Range("L3").Value = "2017-02-21 04:02:14"
$dede.Range("M3").Value = ""
$dede.Range("R3").Value = ""

$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(13).ColumnWidth = 23.0
$dede.Columns.Item(18).ColumnWidth = 12.833333333333332
